$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.677.46'
$ws.Range("E2").Value = '  -2.44%  '
$ws.Range("D3").Value = '3.090.60'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.28'
$ws.Range("E5").Value = '  -3.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.09'
$ws.Range("E6").Value = '  -5.32%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.089.59'
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.145'
$ws.Range("E10").Value = '  -2.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.24'
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("E12").Value = '  -2.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.45'
$ws.Range("E14").Value = '  -4.44%  '
$ws.Range("D15").Value = '3.598.26'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.119'
$ws.Range("D17").Value = '62.694.78'
$ws.Range("E17").Value = '  -2.41%  '
$ws.Range("D18").Value = '3.087.52'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '472.14'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.13'
$ws.Range("E21").Value = '  -2.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.688'
$ws.Range("E22").Value = '  -3.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.55'
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.27'
$ws.Range("E24").Value = '  +3.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.80'
$ws.Range("E25").Value = '  -5.81%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").Value = '  -3.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.96'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.89'
$ws.Range("E29").Value = '  -3.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.01'
$ws.Range("E30").Value = '  -2.71%  '
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.85'
$ws.Range("E32").Value = '  +2.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.106'
$ws.Range("E33").Value = '  -7.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.53'
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.81'
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.81'
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").Value = '0.0₃0711'
$ws.Range("E38").Value = '  -4.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0382'
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '415.66'
$ws.Range("E40").Value = '  -7.20%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  -11.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.15'
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("D44").Value = '2.849.17'
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.253'
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("E47").Value = '  -6.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.29'
$ws.Range("E48").Value = '  -4.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.112'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.07'
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.74'
$ws.Range("E51").Value = '  -0.04%  '
